$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 559
$ws1.Range("F5").Value = 6315
$ws1.Range("F7").Value = 1080
$ws1.Range("F8").Value = 63
$ws1.Range("F9").Value = 618
$ws1.Range("F10").Value = 304
$ws1.Range("F12").Value = 666
$ws1.Range("F13").Value = 2
$ws1.Range("F14").Value = 1132
$ws1.Range("F15").Value = 74
$ws1.Range("F16").Value = 390
$ws1.Range("F18").Value = 14
$ws1.Range("F20").Value = 649
$ws1.Range("F21").Value = 365
$ws1.Range("F22").Value = 385
$ws1.Range("F25").Value = 114
$ws1.Range("F26").Value = 2171
$ws1.Range("F28").Value = 82
$ws1.Range("F31").Value = 3497

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 157
$ws2.Range("F8").Value = 698
$ws2.Range("F20").Value = 4081
$ws2.Range("F24").Value = 180
$ws2.Range("F26").Value = 84
$ws2.Range("F29").Value = 22
$ws2.Range("F32").Value = 1572

# --- Sheet "本地生活" ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("E5").Value = "2024.01.06 00:00-03.31 23:59"
$ws3.Range("F5").Value = 1177
$ws3.Range("F7").Value = 1562
$ws3.Range("F11").Value = 734

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("E4").Value = "2024.01.06 00:00-03.31 23:59"
$ws4.Range("F4").Value = 1177
$ws4.Range("F5").Value = 1562
$ws4.Range("F8").Value = 734
$ws4.Range("F9").Value = 559
$ws4.Range("F12").Value = 6315
$ws4.Range("F16").Value = 1080
$ws4.Range("F17").Value = 698
$ws4.Range("F18").Value = 618
$ws4.Range("F20").Value = 666
$ws4.Range("F25").Value = 1132
$ws4.Range("F26").Value = 390
$ws4.Range("F30").Value = 14
$ws4.Range("F34").Value = 649
$ws4.Range("F35").Value = 365
$ws4.Range("F36").Value = 385
$ws4.Range("F39").Value = 180
$ws4.Range("F41").Value = 84
$ws4.Range("F45").Value = 1572
$ws4.Range("F47").Value = 82
$ws4.Range("F50").Value = 3497
